# Update "想去人数" (want-to-attend count) figures in F column across the
# three sheets that carry them: 展览 (exhibitions), 演出 (performances),
# and 全部类型 (the combined/all-types sheet). 本地生活 is untouched.

$wb = $excel.ActiveWorkbook

# --- 展览 (exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 2952
$ws.Range("F8").Value = 1638
$ws.Range("F11").Value = 107
$ws.Range("F15").Value = 1501
$ws.Range("F16").Value = 6997
$ws.Range("F18").Value = 7162
$ws.Range("F20").Value = 5357
$ws.Range("F21").Value = 3094
$ws.Range("F23").Value = 220
$ws.Range("F24").Value = 166
$ws.Range("F25").Value = 1851
$ws.Range("F26").Value = 76
$ws.Range("F27").Value = 295
$ws.Range("F32").Value = 2391
$ws.Range("F33").Value = 1139
$ws.Range("F34").Value = 2615
$ws.Range("F35").Value = 14
$ws.Range("F39").Value = 1040
$ws.Range("F40").Value = 207

# --- 演出 (performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 209

# --- 全部类型 (combined / all types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 2952
$ws.Range("F9").Value = 1638
$ws.Range("F12").Value = 107
$ws.Range("F15").Value = 1501
$ws.Range("F16").Value = 209
$ws.Range("F19").Value = 6997
$ws.Range("F21").Value = 7162
$ws.Range("F23").Value = 5357
$ws.Range("F24").Value = 3094
$ws.Range("F27").Value = 220
$ws.Range("F29").Value = 1851
$ws.Range("F32").Value = 295
$ws.Range("F37").Value = 2391
$ws.Range("F38").Value = 1139
$ws.Range("F40").Value = 2615
$ws.Range("F41").Value = 14
$ws.Range("F46").Value = 1040
$ws.Range("F47").Value = 207

Write-Output "Applied 38 F-column updates across 展览/演出/全部类型."
